$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Cypher statistics query (shared by C2:C4 "StatQuery" cells) -
# rewritten so the filter counts programs/studies/cases/samples/files
# instead of the old parameterised multi-filter version.
$newQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Golden Retriever']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newQuery
$ws.Range("C3").Value = $newQuery
$ws.Range("C4").Value = $newQuery

# The wrapped query text is now shorter, so the wrapped row heights shrink
# from the capped maximum (409.6) down to their real auto-fit heights.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# Reset the view back to the top-left of the sheet (was scrolled to A4).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
